$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are stored as plain text strings that happen to
# look like numbers, so force text format first to stop Excel from
# re-interpreting them as numeric values on assignment.
$priceChanges = [ordered]@{
    "D2"  = "244.36"
    "D3"  = "21.87"
    "D4"  = "5.459"
    "D5"  = "0.05759"
    "D6"  = "3.417"
    "D7"  = "6.314"
    "D8"  = "0.8142"
    "D9"  = "1.055"
    "D10" = "0.1425"
    "D11" = "0.07284"
    "D12" = "0.03134"
    "D13" = "0.03156"
    "D14" = "4.139"
    "D15" = "0.09373"
    "D16" = "0.001595"
    "D17" = "0.04806"
    "D18" = "0.0005838"
    "D19" = "0.006304"
    "D20" = "0.004135"
    "D21" = "0.0009953"
    "D22" = "0.0001497"
    "D24" = "2.160"
    "D26" = "0.1297"
    "D27" = "0.0003991"
    "D40" = "0.03846"
    "D41" = "0.006656"
    "D42" = "0.1070"
    "D43" = "0.002725"
    "D44" = "0.006557"
    "D46" = "0.00000000749"
    "D47" = "0.3892"
    "D49" = "0.00002096"
}

foreach ($addr in $priceChanges.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $priceChanges[$addr]
}

# Volume(1h) column (E) entries are plain text labels; no numeric coercion
# to worry about, so just assign directly.
$ws.Range("E9").Value = "8FTXTokenFTTBestin24h"
$ws.Range("E18").Value = "17OneONE"
